$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @(65071, "Luiz Gustavo Barros", "Recursos Humanos", "Problemas pessoais", 7, 45096, 3042),
    @(64824, "Ana Júlia Andrade", "Marketing", "Outros", 3, 45104, 2259.55),
    @(52874, "Manuela Lima", "Operacoes", "Consulta medica", 4, 45091, 3651.48),
    @(43213, "Maria Júlia Sá", "Juridico", "Problemas pessoais", 5, 45098, 5707.72),
    @(49712, "Maria Liz Ferreira", "P&D", "Consulta medica", 7, 45087, 6110.01),
    @(74651, "Augusto Casa Grande", "Engenharia", "Consulta medica", 8, 45079, 8154.67),
    @(30076, "Sarah Gonçalves", "Engenharia", "Consulta medica", 4, 45101, 6421.3),
    @(84085, "Ester Siqueira", "Juridico", "Viagem de negocios", 6, 45088, 5896.46),
    @(49142, "Srta. Ágatha da Luz", "Engenharia", "Viagem de negocios", 2, 45093, 8480.92),
    @(79804, "Sr. Diego Sampaio", "Atendimento ao Cliente", "Doenca", 6, 45078, 8412.79)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
